# "Ship via" block: fill in the Shipped Via value for the invoice header.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D12").Value = "test"

# Populate the first (only remaining) line item row with real data.
$ws.Range("B16").Value = "bbbb"
$ws.Range("D16").Value = "3"
$ws.Range("E16").Value = "01008-1404"
$ws.Range("G16").Value = "SADDLE"
$ws.Range("I16").Value = 5000
$ws.Range("J16").Value = "EACH"
$ws.Range("L16").Value = "DI"
$ws.Range("M16").Value = 3.8
$ws.Range("N16").Value = "0%"
$ws.Range("P16").Value = 0
$ws.Range("Q16").Value = 3.8
$ws.Range("R16").Value = 19000

# Remove the second (still-template) line item row entirely - shifts the
# footer rows up by one.
$ws.Range("A17:R17").EntireRow.Delete()

# The "TOTAL DUE" row (now row 17 after the delete) reflects the new total.
$ws.Range("R17").Value = 19000
